# Updates cryptos.xlsx price/volume columns (D, E) and two swapped rows (B, C)
# for rows 2-51, per commit "Updated cryptos list on Wed May 10 08:33:23 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.624.78'
$ws.Range("E2").Value = '  +0.00%  '
# Row 3
$ws.Range("D3").Value = '1.841.75'
$ws.Range("E3").Value = '  -0.14%  '
# Row 4
$ws.Range("E4").Value = '  +0.09%  '
# Row 5
$ws.Range("D5").Value = '''312.16'
$ws.Range("E5").Value = '  -1.17%  '
# Row 6
$ws.Range("D6").Value = '''1.002'
$ws.Range("E6").Value = '  +0.13%  '
# Row 7
$ws.Range("D7").Value = '''0.4255'
# Row 8
$ws.Range("D8").Value = '''0.3609'
$ws.Range("E8").Value = '  -0.87%  '
# Row 9
$ws.Range("D9").Value = '''0.07305'
$ws.Range("E9").Value = '  +0.73%  '
# Row 10
$ws.Range("D10").Value = '''0.8746'
$ws.Range("E10").Value = '  -1.86%  '
# Row 11
$ws.Range("D11").Value = '''20.57'
$ws.Range("E11").Value = '  -0.21%  '
# Row 12
$ws.Range("D12").Value = '1.817.78'
$ws.Range("E12").Value = '  -0.74%  '
# Row 13
$ws.Range("D13").Value = '''5.321'
# Row 14
$ws.Range("D14").Value = '''6.487'
$ws.Range("E14").Value = '  -1.44%  '
# Row 15
$ws.Range("D15").Value = '''0.06969'
$ws.Range("E15").Value = '  +1.44%  '
# Row 16
$ws.Range("E16").Value = '  +0.16%  '
# Row 17
$ws.Range("D17").Value = '''79.24'
$ws.Range("E17").Value = '  +0.72%  '
# Row 18
$ws.Range("D18").Value = '''0.000008927'
$ws.Range("E18").Value = '  +1.07%  '
# Row 19
$ws.Range("D19").Value = '''1.002'
$ws.Range("E19").Value = '  +0.26%  '
# Row 20
$ws.Range("D20").Value = '''15.31'
$ws.Range("E20").Value = '  -0.75%  '
# Row 21
$ws.Range("D21").Value = '27.729.75'
$ws.Range("E21").Value = '  +0.42%  '
# Row 22
$ws.Range("D22").Value = '''4.961'
$ws.Range("E22").Value = '  -0.50%  '
# Row 23
$ws.Range("E23").Value = '  -2.51%  '
# Row 24
$ws.Range("D24").Value = '2.054.75'
$ws.Range("E24").Value = '  +0.29%  '
# Row 25
$ws.Range("D25").Value = '''1.989'
$ws.Range("E25").Value = '  -1.35%  '
# Row 26
$ws.Range("D26").Value = '''155.50'
$ws.Range("E26").Value = '  +0.61%  '
# Row 27
$ws.Range("D27").Value = '''18.50'
$ws.Range("E27").Value = '  -0.34%  '
# Row 28
$ws.Range("D28").Value = '''119.29'
$ws.Range("E28").Value = '  +0.10%  '
# Row 29
$ws.Range("D29").Value = '''5.196'
$ws.Range("E29").Value = '  -0.76%  '
# Row 30
$ws.Range("D30").Value = '''1.873'
$ws.Range("E30").Value = '  +1.64%  '
# Row 31
$ws.Range("D31").Value = '''0.08858'
$ws.Range("E31").Value = '  -0.55%  '
# Row 32
$ws.Range("D32").Value = '''0.7581'
# Row 33
$ws.Range("D33").Value = '''2.954'
$ws.Range("E33").Value = '  -0.03%  '
# Row 34
$ws.Range("D34").Value = '''4.493'
$ws.Range("E34").Value = '  -1.49%  '
# Row 35
$ws.Range("D35").Value = '''1.125'
$ws.Range("E35").Value = '  +2.12%  '
# Row 36
$ws.Range("E36").Value = '  +0.12%  '
# Row 37
$ws.Range("D37").Value = '''0.05416'
$ws.Range("E37").Value = '  +0.38%  '
# Row 38
$ws.Range("D38").Value = '''1.102'
$ws.Range("E38").Value = '  +0.14%  '
# Row 39
$ws.Range("D39").Value = '''0.01923'
$ws.Range("E39").Value = '  -0.07%  '
# Row 40
$ws.Range("D40").Value = '''2.818'
$ws.Range("E40").Value = '  +0.00%  '
# Row 41
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '''0.1656'
$ws.Range("E41").Value = '  +0.43%  '
# Row 42
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '''0.5049'
$ws.Range("E42").Value = '  -0.31%  '
# Row 43
$ws.Range("D43").Value = '''6.530'
$ws.Range("E43").Value = '  -4.95%  '
# Row 44
$ws.Range("D44").Value = '''8.361'
$ws.Range("E44").Value = '  +1.19%  '
# Row 45
$ws.Range("D45").Value = '''0.06548'
$ws.Range("E45").Value = '  -0.89%  '
# Row 46
$ws.Range("D46").Value = '''10.37'
$ws.Range("E46").Value = '  +0.31%  '
# Row 47
$ws.Range("D47").Value = '''105.97'
$ws.Range("E47").Value = '  +1.14%  '
# Row 48
$ws.Range("E48").Value = '  +0.16%  '
# Row 49
$ws.Range("E49").Value = '  -1.70%  '
# Row 50
$ws.Range("D50").Value = '''1.633'
$ws.Range("E50").Value = '  +0.27%  '
# Row 51
$ws.Range("D51").Value = '''64.22'
$ws.Range("E51").Value = '  -0.54%  '
